# Apply the author's edits to the "Isla Natividad Indicators" workbook:
#  1. Highlight (red fill) the indicator-name cells in column B that still
#     need attention/follow-up.
#  2. Leave the selection on B15, matching where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells flagged with a red highlight (new style: solid red fill, no border).
$flaggedCells = @("B5", "B6", "B8", "B9", "B10", "B11", "B13", "B15", "B16", "B17")
foreach ($addr in $flaggedCells) {
    $ws.Range($addr).Interior.Color = 255
}

# Update the active selection/cell to B15.
$ws.Range("B15").Select()
